# Update the "想去人数" (want-to-go count) figures in the F column on the
# sheets that carry the data table: "展览" and "全部类型" (kept in sync).
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 157
    "F3"  = 7101
    "F4"  = 4955
    "F5"  = 72
    "F9"  = 95
    "F10" = 70
    "F11" = 73
    "F13" = 621
    "F14" = 150
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
